# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value pairs scraped from the refreshed coinranking.com feed
$updates = @(
    @{ Cell = 'D2'; Value = '65.802.76' }
    @{ Cell = 'E2'; Value = '  +0.12%  ' }
    @{ Cell = 'D3'; Value = '2.664.43' }
    @{ Cell = 'E3'; Value = '  -0.47%  ' }
    @{ Cell = 'E4'; Value = '  -0.04%  ' }
    @{ Cell = 'D5'; Value = '598.37' }
    @{ Cell = 'E5'; Value = '  -0.36%  ' }
    @{ Cell = 'D6'; Value = '158.06' }
    @{ Cell = 'E6'; Value = '  +0.87%  ' }
    @{ Cell = 'D7'; Value = '0.655' }
    @{ Cell = 'E7'; Value = '  +4.95%  ' }
    @{ Cell = 'E8'; Value = '  -0.03%  ' }
    @{ Cell = 'E9'; Value = '  -2.80%  ' }
    @{ Cell = 'E10'; Value = '  +0.64%  ' }
    @{ Cell = 'D11'; Value = '5.85' }
    @{ Cell = 'E11'; Value = '  -0.32%  ' }
    @{ Cell = 'E12'; Value = '  +1.48%  ' }
    @{ Cell = 'D13'; Value = '28.99' }
    @{ Cell = 'E13'; Value = '  -1.36%  ' }
    @{ Cell = 'E14'; Value = '  -2.44%  ' }
    @{ Cell = 'D15'; Value = '3.141.29' }
    @{ Cell = 'E15'; Value = '  -0.56%  ' }
    @{ Cell = 'D16'; Value = '65.671.50' }
    @{ Cell = 'E16'; Value = '  +0.15%  ' }
    @{ Cell = 'D17'; Value = '2.664.92' }
    @{ Cell = 'E17'; Value = '  -0.46%  ' }
    @{ Cell = 'D18'; Value = '12.58' }
    @{ Cell = 'E18'; Value = '  -2.54%  ' }
    @{ Cell = 'D19'; Value = '4.80' }
    @{ Cell = 'E19'; Value = '  -0.10%  ' }
    @{ Cell = 'B20'; Value = 'Uniswap' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ Cell = 'D20'; Value = '7.49' }
    @{ Cell = 'E20'; Value = '  -1.30%  ' }
    @{ Cell = 'B21'; Value = 'BitcoinCash' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' }
    @{ Cell = 'D21'; Value = '351.44' }
    @{ Cell = 'E21'; Value = '  -0.25%  ' }
    @{ Cell = 'E23'; Value = '  +0.08%  ' }
    @{ Cell = 'D24'; Value = '1.83' }
    @{ Cell = 'E24'; Value = '  +11.88%  ' }
    @{ Cell = 'E25'; Value = '  +0.62%  ' }
    @{ Cell = 'D26'; Value = '9.68' }
    @{ Cell = 'E26'; Value = '  +0.15%  ' }
    @{ Cell = 'E27'; Value = '  +1.75%  ' }
    @{ Cell = 'D28'; Value = '571.57' }
    @{ Cell = 'E28'; Value = '  +7.65%  ' }
    @{ Cell = 'D29'; Value = '8.20' }
    @{ Cell = 'E29'; Value = '  +1.56%  ' }
    @{ Cell = 'E30'; Value = '  -2.87%  ' }
    @{ Cell = 'E31'; Value = '  -0.33%  ' }
    @{ Cell = 'E32'; Value = '  +0.44%  ' }
    @{ Cell = 'E33'; Value = '  +4.34%  ' }
    @{ Cell = 'D34'; Value = '6.70' }
    @{ Cell = 'E34'; Value = '  +3.46%  ' }
    @{ Cell = 'D35'; Value = '5.58' }
    @{ Cell = 'E35'; Value = '  +1.83%  ' }
    @{ Cell = 'E36'; Value = '  -0.34%  ' }
    @{ Cell = 'D37'; Value = '20.64' }
    @{ Cell = 'E37'; Value = '  +0.42%  ' }
    @{ Cell = 'E38'; Value = '  -0.03%  ' }
    @{ Cell = 'E39'; Value = '  +0.11%  ' }
    @{ Cell = 'D40'; Value = '155.16' }
    @{ Cell = 'E40'; Value = '  -2.00%  ' }
    @{ Cell = 'D41'; Value = '161.48' }
    @{ Cell = 'E41'; Value = '  -1.81%  ' }
    @{ Cell = 'E42'; Value = '  -1.11%  ' }
    @{ Cell = 'D43'; Value = '0.0620' }
    @{ Cell = 'E43'; Value = '  +1.68%  ' }
    @{ Cell = 'E44'; Value = '  -0.35%  ' }
    @{ Cell = 'D45'; Value = '23.03' }
    @{ Cell = 'E45'; Value = '  +0.60%  ' }
    @{ Cell = 'D46'; Value = '0.645' }
    @{ Cell = 'E46'; Value = '  +0.25%  ' }
    @{ Cell = 'E47'; Value = '  +0.17%  ' }
    @{ Cell = 'E48'; Value = '  +1.61%  ' }
    @{ Cell = 'D49'; Value = '19.81' }
    @{ Cell = 'E49'; Value = '  -1.68%  ' }
    @{ Cell = 'D50'; Value = '0.0₆0247' }
    @{ Cell = 'E50'; Value = '  -4.38%  ' }
    @{ Cell = 'D51'; Value = '0.816' }
    @{ Cell = 'E51'; Value = '  -0.14%  ' }
)

foreach ($u in $updates) {
    $val = $u.Value
    $range = $ws.Range($u.Cell)
    # Cells in columns D/E hold plain-text numbers/percentages (coinranking
    # renders them pre-formatted, e.g. "65.802.76" / "  +0.12%  "). Excel
    # auto-coerces anything that parses as a number when assigned via .Value,
    # which would corrupt values like "0.0620" or "5.85" (precision/zeros lost).
    # Prefix those with a quote so Excel stores them as text, exactly like typing
    # them in the UI with a leading apostrophe.
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $val
    } else {
        $range.Value = $val
    }
}
